$d = $word.ActiveDocument

# 1. Split "FreeFundraise.com" into "Free" + " " + "Fundraise.com"
$d.Content.Find.Execute("FreeFundraise.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Free Fundraise.com", 2)

# 2. Shrink the 16-space run (right before "Mar 2015") down to 6 spaces.
#    It's the unique 16-space run that is immediately followed by "Mar ".
$d.Content.Find.Execute("                Mar ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "      Mar ", 2)
